# Grade update - midterm 2
#
# 1. Enter the Midterm 2 grades that were missing for Loren Grey (row 6)
#    and Esiete Yismaw Mebratie (row 7).
# 2. Remove Spencer Jenkins (row 10) from the roster entirely, which
#    shifts the remaining students (Penelope Turgen, Kayla Arias) up one
#    row.
# 3. Leave the selection on the last data row, as it was after making
#    these edits.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Newly-graded Midterm 2 scores (column G)
$ws.Range("G6").Formula = "=59/60"
$ws.Range("G7").Formula = "=53/60"

# Spencer Jenkins withdrew / is removed from the gradebook - delete the
# whole row so everyone below moves up.
$ws.Rows.Item(10).Delete()

# Leave the cursor where the roster now ends.
[void]$ws.Range("G11").Select()
